$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# ---------------------------------------------------------------------------
# Insert 5 new rows at row 16 for the components bought on 2012-11-26/30
# (pushes the old spacer rows 16-26 down to 21-31, and the totals row to 33).
# ---------------------------------------------------------------------------
$ws.Range("A16:A20").EntireRow.Insert()

# Copy the A:F cell formatting (alternating stripe style) from the existing
# data rows just above so the new rows match the table's look exactly.
$ws.Range("A13:F13").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A18:F18").PasteSpecial(-4122)
$ws.Range("A20:F20").PasteSpecial(-4122)
$ws.Range("A14:F14").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("A19:F19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A16:F20").RowHeight = 15.75

# Write the new purchases. The product names are entered in the same order
# the brand-new strings first appear in the finished sheet (row 19, then 18,
# then 16) so the shared-string table grows in that order; rows 17 and 20
# reuse already-existing strings.
$ws.Cells.Item(19,2).Value = "Regulador 74RM33"
$ws.Cells.Item(18,2).Value = "SOQUETE torneado 28 pinos"
$ws.Cells.Item(16,2).Value = "74HC125 BUS Line Driver"
$ws.Cells.Item(17,2).Value = "CONECTOR modu 2542"
$ws.Cells.Item(20,2).Value = "TERMINAL para conector modu 22/26AWG"

$ws.Cells.Item(16,1).Value = 41239
$ws.Cells.Item(16,3).Value = 1.5
$ws.Cells.Item(16,4).Value = 3
$ws.Cells.Item(16,5).Value = "Mundial Componentes"
$ws.Cells.Item(16,6).Formula = "=PRODUCT(C16:D16)"

$ws.Cells.Item(17,1).Value = 41243
$ws.Cells.Item(17,3).Value = 0.6
$ws.Cells.Item(17,4).Value = 2
$ws.Cells.Item(17,5).Value = "Mundial Componentes"
$ws.Cells.Item(17,6).Formula = "=PRODUCT(C17:D17)"

$ws.Cells.Item(18,1).Value = 41243
$ws.Cells.Item(18,3).Value = 3
$ws.Cells.Item(18,4).Value = 1
$ws.Cells.Item(18,5).Value = "Mundial Componentes"
$ws.Cells.Item(18,6).Formula = "=PRODUCT(C18:D18)"

$ws.Cells.Item(19,1).Value = 41243
$ws.Cells.Item(19,3).Value = 2.5
$ws.Cells.Item(19,4).Value = 2
$ws.Cells.Item(19,5).Value = "Mundial Componentes"
$ws.Cells.Item(19,6).Formula = "=PRODUCT(C19:D19)"

$ws.Cells.Item(20,1).Value = 41243
$ws.Cells.Item(20,3).Value = 0.1
$ws.Cells.Item(20,4).Value = 16
$ws.Cells.Item(20,5).Value = "Mundial Componentes"
$ws.Cells.Item(20,6).Formula = "=PRODUCT(C20:D20)"

# ---------------------------------------------------------------------------
# The table still only keeps 7 blank spacer rows above the totals line, so
# trim 4 of the now-surplus blanks (old rows 23-26, now shifted to 28-31).
# ---------------------------------------------------------------------------
$ws.Range("A28:A31").EntireRow.Delete()

# Restore the frozen-pane scroll position / active selection for the edited view.
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$ws.Range("C17").Select()
